$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number + report week date range ---
$ws.Range("A8").Value = 'Volume 31   Number  45'
$ws.Range("C9").Value = 'Report Covering the Week  11/4/2024  Through  11/10/2024'

# --- Weekly precinct crime-stat table updates (rows 15-31) ---
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = "'0"
$ws.Range("H15").Value = "***.*"
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = -9.090909090909
$ws.Range("L15").Value = 11.111111111111
$ws.Range("M15").Value = 42.857142857142
$ws.Range("N15").Value = 100
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -26.086956521739
$ws.Range("I16").Value = 209
$ws.Range("J16").Value = 227
$ws.Range("K16").Value = -7.929515418502
$ws.Range("L16").Value = 24.404761904761
$ws.Range("M16").Value = 27.439024390243
$ws.Range("N16").Value = -76.168757126567
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 27.777777777777
$ws.Range("I17").Value = 228
$ws.Range("J17").Value = 213
$ws.Range("K17").Value = 7.042253521126
$ws.Range("L17").Value = 17.525773195876
$ws.Range("M17").Value = 107.272727272727
$ws.Range("N17").Value = -15.555555555555
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -38.095238095238
$ws.Range("I18").Value = 210
$ws.Range("J18").Value = 187
$ws.Range("K18").Value = 12.299465240641
$ws.Range("L18").Value = 43.835616438356
$ws.Range("M18").Value = -1.408450704225
$ws.Range("N18").Value = -84.126984126984
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 77.777777777777
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 10.204081632653
$ws.Range("I19").Value = 628
$ws.Range("J19").Value = 643
$ws.Range("K19").Value = -2.332814930015
$ws.Range("L19").Value = 3.630363036303
$ws.Range("M19").Value = 55.831265508684
$ws.Range("N19").Value = -20.807061790668
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -34.482758620689
$ws.Range("I20").Value = 225
$ws.Range("J20").Value = 270
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = 10.294117647058
$ws.Range("M20").Value = 24.309392265193
$ws.Range("N20").Value = -87.596471885336
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 37.037037037037
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 140
$ws.Range("H21").Value = -6.428571428571
$ws.Range("I21").Value = 1521
$ws.Range("J21").Value = 1562
$ws.Range("K21").Value = -2.624839948783
$ws.Range("L21").Value = 13.762154076290
$ws.Range("M21").Value = 39.926402943882
$ws.Range("N21").Value = -70.147203140333
$ws.Range("F22").Value = 8
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 53
$ws.Range("J22").Value = 75
$ws.Range("K22").Value = -29.333333333333
$ws.Range("L22").Value = -26.388888888888
$ws.Range("M22").Value = 20.454545454545
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = -39.130434782608
$ws.Range("F24").Value = 137
$ws.Range("G24").Value = 193
$ws.Range("H24").Value = -29.015544041450
$ws.Range("I24").Value = 1772
$ws.Range("J24").Value = 1778
$ws.Range("K24").Value = -0.337457817772
$ws.Range("L24").Value = 35.993860322333
$ws.Range("M24").Value = 116.361416361416
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 89
$ws.Range("G25").Value = 121
$ws.Range("H25").Value = -26.446280991735
$ws.Range("I25").Value = 1196
$ws.Range("J25").Value = 1104
$ws.Range("K25").Value = 8.333333333333
$ws.Range("L25").Value = 123.551401869159
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -36.363636363636
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 28.205128205128
$ws.Range("I26").Value = 487
$ws.Range("J26").Value = 448
$ws.Range("K26").Value = 8.705357142857
$ws.Range("L26").Value = -0.612244897959
$ws.Range("M26").Value = 11.954022988505
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 29
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 20.833333333333
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -28.571428571428
$ws.Range("I28").Value = 59
$ws.Range("J28").Value = 75
$ws.Range("K28").Value = -21.333333333333
$ws.Range("L28").Value = -26.25
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").Value = "'0"
$ws.Range("H31").Value = -100

Write-Output "edits applied"
